$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 5): roster names ---
# D5: name1 (leader) -> Ryan Conyac
$ws.Range("D5").Value = "Ryan Conyac"
# E5 (Michael McGregor) and F5 (Yuchen Feng) stay the same
# G5: name4 -> Will McLain
$ws.Range("G5").Value = "Will McLain"
# H5: name5 -> Younouss Thiam
$ws.Range("H5").Value = "Younouss Thiam"
# I5: name6 -> Brian Davis
$ws.Range("I5").Value = "Brian Davis"
# J5 / K5: name7 / name8 removed entirely (team is down to 6 members)
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""

# --- Row 6 (first meeting) ---
$ws.Range("B6").Value = "8/22 / 1:00"
$ws.Range("D6").Value = "A"
$ws.Range("G6").Value = "A"
$ws.Range("H6").Value = "E"
$ws.Range("I6").Value = "A"

# --- Row 7 (second meeting) ---
$ws.Range("B7").Value = "8/25 / 4:15"
$ws.Range("D7").Value = "A"
$ws.Range("G7").Value = "A"
$ws.Range("H7").Value = "A"
$ws.Range("I7").Value = "E"

# --- Row 8 / Row 9 (upcoming meeting dates) ---
$ws.Range("B8").Value = "8/29 / 1:00"
$ws.Range("B9").Value = "9/1 / 4:15"

# Move the active selection to where the author left off editing
$ws.Range("B9").Select()
